# Regenerate save_data: write new "K" (strikeout) values into column G,
# replacing the old Strike# values, for rows 2-39 of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 1
    6  = 1
    7  = 2
    8  = 1
    9  = 3
    10 = 2
    11 = 1
    12 = 1
    13 = 1
    14 = 3
    15 = 2
    16 = 2
    17 = 0
    18 = 2
    19 = 3
    20 = 1
    21 = 0
    22 = 0
    23 = 1
    24 = 6
    25 = 2
    26 = 1
    27 = 3
    28 = 3
    29 = 2
    30 = 2
    31 = 8
    32 = 3
    33 = 1
    34 = 3
    35 = 3
    36 = 0
    37 = 1
    38 = 2
    39 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
